# Auto-generated script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.005.27"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "2.648.37"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'569.54"
$ws.Range("E5").Value = "  +6.48%  "
$ws.Range("D6").Value = "'147.18"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  +7.37%  "
$ws.Range("D9").Value = "2.674.61"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "'6.85"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +4.62%  "
$ws.Range("E12").Value = "  +6.42%  "
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "3.117.82"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "60.883.08"
$ws.Range("E15").Value = "  +4.15%  "
$ws.Range("D16").Value = "'21.89"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +4.61%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.667.50"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "'344.76"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("D21").Value = "'10.48"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "'6.39"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'5.82"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'66.90"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'0.441"
$ws.Range("E26").Value = "  +6.30%  "
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("D30").Value = "0.0₃0786"
$ws.Range("E30").Value = "  +7.02%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("D33").Value = "'6.28"
$ws.Range("E33").Value = "  +7.14%  "
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "'154.68"
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").Value = "'4.10"
$ws.Range("E36").Value = "  +5.43%  "
$ws.Range("E37").Value = "  +8.10%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.907"
$ws.Range("E38").Value = "  +6.93%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.912"
$ws.Range("E39").Value = "  +12.23%  "
$ws.Range("D40").Value = "'37.59"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  +7.45%  "
$ws.Range("D42").Value = "'303.91"
$ws.Range("E42").Value = "  +7.88%  "
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.608"
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0984"
$ws.Range("E46").Value = "  +5.02%  "
$ws.Range("D47").Value = "'0.0551"
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("D48").Value = "'129.02"
$ws.Range("E48").Value = "  +13.27%  "
$ws.Range("D49").Value = "'19.55"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").Value = "'10.66"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  +4.93%  "
